$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.326.76"
$ws.Range("E2").Value = "  -1.54%  "
$ws.Range("D3").Value = "2.634.24"
$ws.Range("E3").Value = "  -2.84%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.12"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.50"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("E8").Value = "  -3.79%  "
$ws.Range("D9").Value = "2.634.17"
$ws.Range("E10").Value = "  -1.50%  "
$ws.Range("E11").Value = "  +1.25%  "
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.26"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.84"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.90%  "
$ws.Range("D15").Value = "3.114.20"
$ws.Range("E15").Value = "  -2.86%  "
$ws.Range("E16").Value = "  -2.79%  "
$ws.Range("D17").Value = "67.295.17"
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("D18").Value = "2.620.38"
$ws.Range("E18").Value = "  -3.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.17"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.06"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +5.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "360.65"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.36"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.70"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.90"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +9.20%  "
$ws.Range("E25").Value = "  -5.30%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "70.45"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.48%  "
$ws.Range("E28").Value = "  -2.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  -2.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "556.74"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.94"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.05%  "
$ws.Range("E33").Value = "  -2.93%  "
$ws.Range("E34").Value = "  -3.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.136"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.01%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  -4.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.58"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.20"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.29%  "
$ws.Range("E40").Value = "  -2.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.22"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.88%  "
$ws.Range("E42").Value = "  -3.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.94"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  -4.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.15"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.51%  "
$ws.Range("D47").Value = "0.0₆0301"
$ws.Range("E47").Value = "  -2.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.588"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "152.16"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.66%  "
$ws.Range("E50").Value = "  -1.88%  "
$ws.Range("E51").Value = "  -1.55%  "
